$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 43

# Columns A-D hold text (date/time/weekday/week-number strings). Force
# text formatting first so Excel doesn't auto-coerce these into a date
# serial number / time fraction / plain integer, matching every other
# data row already in the sheet.
$textRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 4))
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-01-26"
$ws.Cells.Item($row, 2).Value = "18:37:18"
$ws.Cells.Item($row, 3).Value = "Sunday"
$ws.Cells.Item($row, 4).Value = "04"

# Restore the default (unstyled) cell style now that the values have
# been committed as text, so the new row matches the plain look of the
# other data rows.
$textRange.Style = "Normal"

$ws.Cells.Item($row, 5).Value = 126108
$ws.Cells.Item($row, 6).Value = 142021
$ws.Cells.Item($row, 7).Value = 168036
$ws.Cells.Item($row, 8).Value = 158489
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142567
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191541
$ws.Cells.Item($row, 14).Value = 115618
$ws.Cells.Item($row, 15).Value = 45499
$ws.Cells.Item($row, 16).Value = 28373
$ws.Cells.Item($row, 17).Value = 64907
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 45888
$ws.Cells.Item($row, 20).Value = -1
